$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns for season record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/bordered/centered header style used by the rest of row 1 (e.g. A1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record (Wins, Losses, Ties) is constant for every player row on this sheet
$wins = 97
$losses = 65
$ties = 0

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
